# Apply the edits described by the diff:
#  1. Replace the computed "K" column formulas (=J-I, with row 4 holding
#     the shared-formula anchor for K4:K22) with hard-coded literal
#     values - i.e. the formulas were "pasted as values".
#  2. Update the sheet view: scroll so column L is the left-most visible
#     column, and change the selection from U1:V1048576 to M1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# --- 1. Replace formulas in K2:K22 with literal values -----------------
$kValues = @(
    @(2,  4.5499999999999972),
    @(3,  11.72999999999999),
    @(4,  29.36),
    @(5,  11.489999999999995),
    @(6,  11.969999999999999),
    @(7,  24.730000000000004),
    @(8,  22.569999999999993),
    @(9,  11.829999999999998),
    @(10, 14.739999999999995),
    @(11, -2.6899999999999977),
    @(12, 12.799999999999997),
    @(13, 18.11),
    @(14, 20.93),
    @(15, 13.300000000000011),
    @(16, 19.190000000000012),
    @(17, 26.729999999999997),
    @(18, 11.979999999999997),
    @(19, 17.159999999999997),
    @(20, 14.350000000000009),
    @(21, 15.719999999999999),
    @(22, 19.530000000000008)
)

foreach ($pair in $kValues) {
    $row = $pair[0]
    $val = $pair[1]
    $ws.Cells.Item($row, 11).Value = $val
}

# --- 2. Update the sheet view / selection -------------------------------
$ws.Range("M1").Select()
$win = $excel.ActiveWindow
$win.ScrollColumn = 12
$win.ScrollRow = 1
